$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the MODS wrapper element: the template previously wrapped the
# descriptive metadata in <update type="MODS">...</update>; it now uses
# <datastream type="md_descriptive" operation="update">...</datastream>.
$ws.Range("C2").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">'
$ws.Range("Y2").Value = '</mods:mods></datastream></object>'
